# "Generate Report for Handback"
#
# The localization-status report is regenerated after a successful
# handback: the per-language status moves from "Ready for handoff" to
# "Handed back: in sync with en-US", the Latest Handback DateTime for
# each language is refreshed, and the stale "handback file is not the
# latest" Error Detail message is cleared now that the handback is current.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet: language status cells ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = "Handed back: in sync with en-US"
$ov.Range("F2").Value = "Handed back: in sync with en-US"

# Columns widen to fit the new, longer status text. The target widths
# (29.9777047293527 / 13.7470528738839 "characters") come from the source
# report generator's own fit metric and don't sit on Excel's internal
# 1/6-character column-width grid, so we dial in the nearest value that
# rounds, through that grid, back to the closest possible width.
$ov.Columns.Item(5).ColumnWidth = 29.1666666666667
$ov.Columns.Item(6).ColumnWidth = 29.1666666666667

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("K2").Value = "2016-09-05 04:55:16"
$zh.Range("P2").Value = ""

$zh.Columns.Item(3).ColumnWidth = 29.1666666666667
$zh.Columns.Item(16).ColumnWidth = 12.8333333333333

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("K2").Value = "2016-09-05 04:55:24"
$de.Range("P2").Value = ""

$de.Columns.Item(3).ColumnWidth = 29.1666666666667
$de.Columns.Item(16).ColumnWidth = 12.8333333333333
